# SU plot bump  LAS
# Update the "bump" sheet (sheet1) with the refreshed LAS sweep data and
# move the flattened plot-source block from rows 22:34 up to rows 18:30.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Update the raw sweep results (columns E/F) for rows 2-14 -----------
$ws.Range("E2").Value  = 237.18
$ws.Range("F2").Value  = 452.29
$ws.Range("E3").Value  = 238.46
$ws.Range("F3").Value  = 449.01
$ws.Range("E4").Value  = 239.51
$ws.Range("F4").Value  = 445.68
$ws.Range("E5").Value  = 240.32
$ws.Range("F5").Value  = 442.31
$ws.Range("E6").Value  = 240.9
$ws.Range("F6").Value  = 438.9
$ws.Range("E7").Value  = 241.25
$ws.Range("F7").Value  = 435.46
$ws.Range("E8").Value  = 241.36
$ws.Range("F8").Value  = 432
$ws.Range("E9").Value  = 241.24
$ws.Range("F9").Value  = 428.52
$ws.Range("E10").Value = 240.89
$ws.Range("F10").Value = 425.01
$ws.Range("E11").Value = 240.28
$ws.Range("F11").Value = 421.49
$ws.Range("E12").Value = 239.4
$ws.Range("F12").Value = 417.95
$ws.Range("E13").Value = 238.25
$ws.Range("F13").Value = 414.39
$ws.Range("E14").Value = 236.79
$ws.Range("F14").Value = 410.81

# --- 2. Re-establish I2:I14 / J2:J14 as shared formulas ---------------------
# (setting the formula on the whole block makes the engine emit a shared
# formula group, matching how Excel compacts repeated formulas)
$ws.Range("I2:I14").Formula = "=F2-`$F`$8"
$ws.Range("J2:J14").Formula = "=(I2-I3)/(A2-A3)"

# --- 3. Collapse the blank gap (rows 15-21) so the pasted-values plot block
#        that used to live at rows 22:34 shifts up to rows 18:30 -----------
$ws.Range("A15:A18").EntireRow.Delete() | Out-Null

# --- 4. Refresh the pasted-values block (I18:J30) with the new numbers -----
$ws.Range("I18").Value = 237.18
$ws.Range("J18").Value = 452.29
$ws.Range("I19").Value = 238.46
$ws.Range("J19").Value = 449.01
$ws.Range("I20").Value = 239.51
$ws.Range("J20").Value = 445.68
$ws.Range("I21").Value = 240.32
$ws.Range("J21").Value = 442.31
$ws.Range("I22").Value = 240.9
$ws.Range("J22").Value = 438.9
$ws.Range("I23").Value = 241.25
$ws.Range("J23").Value = 435.46
$ws.Range("I24").Value = 241.36
$ws.Range("J24").Value = 432
$ws.Range("I25").Value = 241.24
$ws.Range("J25").Value = 428.52
$ws.Range("I26").Value = 240.89
$ws.Range("J26").Value = 425.01
$ws.Range("I27").Value = 240.28
$ws.Range("J27").Value = 421.49
$ws.Range("I28").Value = 239.4
$ws.Range("J28").Value = 417.95
$ws.Range("I29").Value = 238.25
$ws.Range("J29").Value = 414.39
$ws.Range("I30").Value = 236.79
$ws.Range("J30").Value = 410.81

# --- 5. Re-anchor the sortState on the moved block (E18:J30, sorted by E
#        descending) so the worksheet's recorded sort matches its new home -
$sortObj = $ws.Sort
$sortObj.SortFields.Clear() | Out-Null
$sortObj.SortFields.Add($ws.Range("E18"), 0, 2) | Out-Null
$sortObj.SetRange($ws.Range("E18:J30"))
$sortObj.Header = 2
$sortObj.Apply()
